$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "1_Vocab_Ex" sheet: selection becomes the whole used range (A1:D7)
# ---------------------------------------------------------------------------
$wsVocabEx = $wb.Worksheets.Item("1_Vocab_Ex")
$wsVocabEx.Range("A1:D7").Select()

# ---------------------------------------------------------------------------
# 2. "2_Vocab" sheet: no longer the active tab, selection becomes A1:B7
# ---------------------------------------------------------------------------
$wsVocab2 = $wb.Worksheets.Item("2_Vocab")
$wsVocab2.Range("A1:B7").Select()

# ---------------------------------------------------------------------------
# 3. "3_" sheet: used to be empty, now gets the "code grammar" matching
#    question.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("3_")

$ws3.Range("A1:B6").WrapText = $true
$ws3.Range("C2:C6").WrapText = $true
$ws3.Range("B7").WrapText = $true

$ws3.Range("A1").Value = "Match the code grammar for a function with what it does that grammar indicates"

$ws3.Range("A2").Value = "Indicates the ""input"" of the function"
$ws3.Range("B2").Value = "C"
$ws3.Range("C2").Value = """return"""

$ws3.Range("A3").Value = "Indicates that the output of the function to the right is being ""assigned"" to the variable on the left"
$ws3.Range("B3").Value = "B"
$ws3.Range("C3").Value = """def"""

$ws3.Range("A4").Value = "Indicates the lines that define the what the function does to the input"
$ws3.Range("B4").Value = "D"
$ws3.Range("C4").Value = "Parentheses"

$ws3.Range("A5").Value = "Indicates the ""output"" of the function (and the end of the function definition)"
$ws3.Range("B5").Value = "C"
$ws3.Range("C5").Value = """="" (equal sign)"

$ws3.Range("A6").Value = "Indicates the beginning of a function definition"
$ws3.Range("B6").Value = "A"
$ws3.Range("C6").Value = "Indented lines in function definition"

$ws3.Range("B7").Value = "B"

$ws3.Rows.Item(1).RowHeight = 60
$ws3.Rows.Item(2).RowHeight = 30
$ws3.Rows.Item(3).RowHeight = 60
$ws3.Rows.Item(4).RowHeight = 45
$ws3.Rows.Item(5).RowHeight = 45
$ws3.Rows.Item(6).RowHeight = 30

$ws3.Range("C7").Select()

# ---------------------------------------------------------------------------
# 4. New sheet "4_" inserted right after "3_" and before "4_MultC" - the
#    "which lines show a function being called" matching/select question.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "4_"

$ws4.Range("A1:B4").WrapText = $true
$ws4.Range("A5:C6").WrapText = $true
$ws4.Range("A7:B7").WrapText = $true

$ws4.Range("A1").Value = "Which of the following lines show a function being ""called""?  Answer all that are true."

$ws4.Range("A2").Value = "y = np.sin(5)"
$ws4.Range("B2").Value = "Y"

$ws4.Range("A3").Value = "def this_is_a_function(time):"
$ws4.Range("B3").Value = "N"

$ws4.Range("A4").Value = "t = t_1 + t_2"
$ws4.Range("B4").Value = "N"

$ws4.Range("A5").Value = "gorilla = elephant(tiger)"
$ws4.Range("B5").Value = "Y"
$ws4.Range("C5").Value = "The grammar is all that matters.  Even though this doesn't seem to make sense, the parentheses indicate that a function is being called and the output is being assigned to the variable ""gorilla"""

$ws4.Range("A6").Value = "money = this_is_a_function(time)"
$ws4.Range("B6").Value = "Y"
$ws4.Range("C6").Value = "There is no ""def"" here, so the parentheses indicate that a function is being called."

$ws4.Rows.Item(1).RowHeight = 45
$ws4.Rows.Item(5).RowHeight = 90
$ws4.Rows.Item(6).RowHeight = 45

$ws4.Columns.Item(1).ColumnWidth = 37.7
$ws4.Columns.Item(3).ColumnWidth = 34.25

# ---------------------------------------------------------------------------
# 5. New blank sheet "Sheet1" inserted right after "4_".
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Add($null, $ws4)
$ws5.Name = "Sheet1"

# ---------------------------------------------------------------------------
# 6. "4_" becomes the active sheet / tab (activeTab moves from 2 to 4).
# ---------------------------------------------------------------------------
$ws4.Range("C6").Select()
